$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear D2, E2, F2; set G2/H2/I2
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "11"
$ws.Range("H2").Value = "'06"
$ws.Range("I2").Value = "1991"

# Row 3: clear D3, E3, F3; G3/H3/I3 unchanged
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4: clear D4, E4, F4; update G4/H4/I4
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "10"
$ws.Range("H4").Value = "81"
$ws.Range("I4").Value = "1990"

# Row 5: clear D5, E5; F5 already empty; G5/H5/I5 unchanged
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
